$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.01124924817744999
$ws.Range("C2").Value = 0.3250033953551542
$ws.Range("D2").Value = 0.4690748481128704
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.2313311609357522
$ws.Range("I2").Value = 0.6767049423462985

# Row 3
$ws.Range("B3").Value = 0.3562861756232065
$ws.Range("C3").Value = 0.2818984299981965
$ws.Range("D3").Value = 0.488069352733541
$ws.Range("E3").Value = 0.05339431913350126
$ws.Range("F3").Value = 0.04774712108833655
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.4839490503001195
$ws.Range("I3").Value = 0.1395325829597731

# Row 4
$ws.Range("B4").Value = 0.1845456900318564
$ws.Range("C4").Value = 0.1321862666885196
$ws.Range("D4").Value = 0.1699182720119316
$ws.Range("E4").Value = 0.151613317885106
$ws.Range("F4").Value = 0.2321718655443033
$ws.Range("G4").Value = 0.1608706819801428
$ws.Range("H4").Value = 0.1387576919398306
$ws.Range("I4").Value = 0.4038575535856002

# Row 5
$ws.Range("B5").Value = 0.3085320317743733
$ws.Range("C5").Value = 0.8491380323169346
$ws.Range("D5").Value = 0.6538963727547132
$ws.Range("E5").Value = -0.02441110408077751
$ws.Range("F5").Value = 0.06377086857981945
$ws.Range("G5").Value = 0.06377086857981945
$ws.Range("H5").Value = 0.3680954120143309
$ws.Range("I5").Value = -0.01901695892664874

# Row 6
$ws.Range("B6").Value = 0.3678593848085374
$ws.Range("C6").Value = 0.2047576031315868
$ws.Range("D6").Value = 0.3785868559086702
$ws.Range("E6").Value = 0.1108905275099262
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.5589600742804085
$ws.Range("I6").Value = 0.3838143409699494

# Row 7
$ws.Range("B7").Value = 0.009318380708520969
$ws.Range("C7").Value = 0.1552079724667571
$ws.Range("D7").Value = 0.01038648203248693
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = -0.01280175566934875
$ws.Range("H7").Value = 0.118959523504126
$ws.Range("I7").Value = -0.01280175566934875

# Row 8
$ws.Range("B8").Value = 0.09996073751251677
$ws.Range("C8").Value = 0.1768414935347959
$ws.Range("D8").Value = 0.2756025859783603
$ws.Range("E8").Value = 0.0225238498398859
$ws.Range("F8").Value = -0.02486276456027156
$ws.Range("G8").Value = -0.01058170630014085
$ws.Range("H8").Value = 0.1492403433014898
$ws.Range("I8").Value = 0.09707718452990623

# Row 9
$ws.Range("B9").Value = 0.296685726588274
$ws.Range("C9").Value = 0.0693593639405216
$ws.Range("D9").Value = 0.01570140263607593
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = -0.008735531775496833
$ws.Range("H9").Value = 0.2091178390319828
$ws.Range("I9").Value = 0.8012695630330529

# Row 10
$ws.Range("B10").Value = -0.01587301587301587
$ws.Range("C10").Value = -0.01587301587301587
$ws.Range("D10").Value = -0.00678179402922467
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0.1321566939302625
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.09954158480681068
$ws.Range("I10").Value = 0.2043222003929273

# Row 11
$ws.Range("B11").Value = 0.008934597828698257
$ws.Range("C11").Value = 0.1294213225719503
$ws.Range("D11").Value = 0.07028549104016626
$ws.Range("E11").Value = 0.09535205945333163
$ws.Range("F11").Value = 0.2423727408389985
$ws.Range("G11").Value = 0.4915423685964469
$ws.Range("H11").Value = 0.126038581154617
$ws.Range("I11").Value = 0.155808664068974

# Row 12
$ws.Range("B12").Value = -0.03802281368821287
$ws.Range("C12").Value = 0.08950874271440465
$ws.Range("D12").Value = -0.02574926129168418
$ws.Range("E12").Value = 0.05204460966542741
$ws.Range("F12").Value = 0.005728314238952386
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = -0.03802281368821287
$ws.Range("I12").Value = 0.02302631578947371

# Row 13
$ws.Range("B13").Value = -0.04904632152588552
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0.01621621621621616
$ws.Range("E13").Value = 0.3514986376021799
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = -0.04904632152588552
$ws.Range("I13").Value = -0.04297994269340966

# Row 14
$ws.Range("B14").Value = -0.04724964739069114
$ws.Range("C14").Value = 0.04281767955801098
$ws.Range("D14").Value = -0.0246575342465754
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = -0.01035911602209953
$ws.Range("I14").Value = -0.06480000000000009

# Row 15
$ws.Range("B15").Value = 0.3966224366706876
$ws.Range("C15").Value = 0.4608996539792387
$ws.Range("D15").Value = 0.4416555407209614
$ws.Range("E15").Value = 0.1496892378148512
$ws.Range("F15").Value = 0.05730791450028883
$ws.Range("G15").Value = 0.02330398757120661
$ws.Range("H15").Value = 0.4436125409452503
$ws.Range("I15").Value = 0.03703882324571978

# Row 16
$ws.Range("B16").Value = 0.3172122627343014
$ws.Range("C16").Value = 0.4383698537693369
$ws.Range("D16").Value = 0.2769765285815673
$ws.Range("E16").Value = 0.01685058787604655
$ws.Range("F16").Value = 0.05811732497915641
$ws.Range("G16").Value = 0.05869619161335039
$ws.Range("H16").Value = 0.3124306021880408
$ws.Range("I16").Value = 0.0004004700809494364

# Row 17
$ws.Range("B17").Value = -0.01565806178586546
$ws.Range("C17").Value = 0.07204780742361756
$ws.Range("D17").Value = 0.03059934249346705
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0.1720430107526882
$ws.Range("I17").Value = -0.0461798583958576

# Row 18
$ws.Range("B18").Value = 0.02572500082346358
$ws.Range("C18").Value = 0.3346942667647007
$ws.Range("D18").Value = 0.03489933096170147
$ws.Range("E18").Value = 0.05224702206900309
$ws.Range("F18").Value = 0.06587493922687847
$ws.Range("G18").Value = 0.02248230549164055
$ws.Range("H18").Value = 0.06483324433984185
$ws.Range("I18").Value = 0.012486023980277

# Row 19
$ws.Range("B19").Value = 0.3623591442757568
$ws.Range("C19").Value = 0.4905143037963895
$ws.Range("D19").Value = 0.5127016997921923
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0.4789222392746438
$ws.Range("I19").Value = 0.2368358824679002

# Row 20
$ws.Range("B20").Value = 0.196566104369388
$ws.Range("C20").Value = 0.4327097163548582
$ws.Range("D20").Value = 0.5530566037735849
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = -0.005432937181663759
$ws.Range("I20").Value = -0.01350084380273768

# Row 21
$ws.Range("B21").Value = 0.007076972566241697
$ws.Range("C21").Value = -0.0939000102134613
$ws.Range("D21").Value = 0.007904385353419764
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0.02261762233850496
$ws.Range("I21").Value = -0.006359522004726558

# Row 22
$ws.Range("B22").Value = 0.05843062175308097
$ws.Range("C22").Value = 0.03197361422437982
$ws.Range("D22").Value = 0.09803297886105697
$ws.Range("E22").Value = -0.01099419559389791
$ws.Range("F22").Value = 0.003717371771517567
$ws.Range("G22").Value = 0.002939416839041448
$ws.Range("H22").Value = 0.02662584223267401
$ws.Range("I22").Value = -0.02695256059862047

